$wb = $excel.ActiveWorkbook

# ============================================================
# Sheet "Erläuterungen" (sheet1) - update explanatory text
# ============================================================
$wsErl = $wb.Worksheets.Item("Erläuterungen")
$wsErl.Range("A4").Value2 = @'
Bisher haben sich mehr als 250 Labore für die RKI-Testlaborabfrage oder in einem der anderen übermittelnden Netzwerke registriert und übermitteln nach Aufruf überwiegend wöchentlich. Da Labore in der RKI-Testzahlabfrage  die  Tests  der  vergangenen Kalenderwochen nachmelden bzw. korrigieren können, ist es möglich, dass sich die ermittelten Zahlen nachträglich ändern. Es ist zu beachten, dass die Zahl der Tests nicht mit der Zahl der getesteten Personen gleichzusetzen ist, da in den Angaben Mehrfachtestungen von Patienten enthalten sein können (s. Testzahlen).
'@
# A5 keeps its existing text ("Zusätzlich zur Anzahl...") - shared string index
# shift happens automatically since the old A4 string is no longer referenced
# anywhere else in the workbook once we overwrite A4.

# View: Excel now shows row 4 at top with A4 selected
$wsErl.Activate() | Out-Null
$wsErl.Range("A4").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 4

# ============================================================
# Sheet "Testzahlen" (sheet2) - weekly test numbers
# ============================================================
$wsTZ = $wb.Worksheets.Item("Testzahlen")

# --- Corrected figures for KW30 (row 24) ---
$wsTZ.Range("C24").Value2 = 553429
$wsTZ.Range("D24").Value2 = 4458
$wsTZ.Range("F24").Value2 = 182

# --- Corrected figures for KW43-KW47 (rows 37-41); values pasted (not formulas) ---
$wsTZ.Range("C37").Value2 = 1418726
$wsTZ.Range("D37").Value2 = 78106
$wsTZ.Range("E37").Value2 = 5.5053618528172459
$wsTZ.Range("F37").Value2 = 204

$wsTZ.Range("C38").Value2 = 1631343
$wsTZ.Range("D38").Value2 = 116673
$wsTZ.Range("E38").Value2 = 7.1519600721614038
$wsTZ.Range("F38").Value2 = 204

$wsTZ.Range("C39").Value2 = 1608125
$wsTZ.Range("D39").Value2 = 126141
$wsTZ.Range("E39").Value2 = 7.8439797901282553
$wsTZ.Range("F39").Value2 = 203

$wsTZ.Range("C40").Value2 = 1396088
$wsTZ.Range("D40").Value2 = 125200
$wsTZ.Range("E40").Value2 = 8.967916062597773
$wsTZ.Range("F40").Value2 = 199

$wsTZ.Range("C41").Value2 = 1363701
$wsTZ.Range("D41").Value2 = 127330
$wsTZ.Range("E41").Value2 = 9.3370907552315341
$wsTZ.Range("F41").Value2 = 198

# --- Insert two new weeks (KW48*, KW49*) before the "Summe" row ---
$wsTZ.Rows.Item(42).Insert() | Out-Null
$wsTZ.Rows.Item(42).Insert() | Out-Null

# Copy the row-41 formatting down into the two freshly inserted rows
$wsTZ.Range("B41:F41").Copy() | Out-Null
$wsTZ.Range("B42:F43").PasteSpecial(-4122) | Out-Null

$wsTZ.Range("B42").Value2 = "48*"
$wsTZ.Range("C42").Value2 = 1340093
$wsTZ.Range("D42").Value2 = 124511
$wsTZ.Range("E42").Value2 = 9.2912208331809811
$wsTZ.Range("F42").Value2 = 199

$wsTZ.Range("B43").Value2 = "49*"
$wsTZ.Range("C43").Value2 = 1297303
$wsTZ.Range("D43").Value2 = 132961
$wsTZ.Range("E43").Value2 = 10.2490320302967
$wsTZ.Range("F43").Value2 = 193

# --- "Summe" row (now row 44): totals pasted as values ---
$wsTZ.Range("C44").Value2 = 30494036
$wsTZ.Range("D44").Value2 = 1252323

# --- Footnote row (now row 45): updated wording ---
$wsTZ.Range("B45").Value2 = @'
*Ab 03. November 2020 geänderte Testkriterien, Daten nicht direkt mit Vorwochen vergleichbar
'@

$wsTZ.Range("B30").Select() | Out-Null
$wsTZ.Range("I42").Select() | Out-Null

# ============================================================
# Sheet "Testkapazitäten" (sheet3) - capacity survey, add KW49 & KW50
# ============================================================
$wsTK = $wb.Worksheets.Item("Testkapazitäten")

$wsTK.Range("A40:E40").Copy() | Out-Null
$wsTK.Range("A41:E42").PasteSpecial(-4122) | Out-Null

$wsTK.Range("A41").Value2 = "KW49"
$wsTK.Range("B41").Value2 = 168
$wsTK.Range("C41").Value2 = 318746
$wsTK.Range("D41").Value2 = 2136828
$wsTK.Range("E41").Value2 = 1918794

$wsTK.Range("A42").Value2 = "KW50"
$wsTK.Range("B42").Value2 = 175
$wsTK.Range("C42").Value2 = 331036
$wsTK.Range("D42").Value2 = 2219158
$wsTK.Range("E42").Value2 = 1944190
# C42/D42 use the same visual style as column E (s=33) in the source workbook
$wsTK.Range("C42:D42").NumberFormat = $wsTK.Range("E42").NumberFormat

$wsTK.Range("A26").Select() | Out-Null
$wsTK.Range("E42").Select() | Out-Null

# ============================================================
# Sheet "Probenrückstau" (sheet4) - backlog survey, add KW48 & KW49
# ============================================================
$wsPR = $wb.Worksheets.Item("Probenrückstau")

$wsPR.Range("A34:C34").Copy() | Out-Null
$wsPR.Range("A35:C36").PasteSpecial(-4122) | Out-Null

$wsPR.Range("A35").Value2 = 53
$wsPR.Range("B35").Value2 = 48
$wsPR.Range("C35").Value2 = 14091

$wsPR.Range("A36").Value2 = 50
$wsPR.Range("B36").Value2 = 49
$wsPR.Range("C36").Value2 = 12237

$wsPR.Range("F42").Select() | Out-Null
